$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# --- G1: new header cell "PRESUPUESTO", cloning the formatting of F1 (bold, bordered, centered header style) ---
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "PRESUPUESTO"

# --- G2: new data cell (0), cloning the formatting of F2 (currency number format) ---
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = 0

# --- G3: new data cell (0), cloning the formatting of F3 (currency, right aligned) ---
$ws.Range("F3").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G3").Value = 0

# --- Column G width = 17 (ColumnWidth setter adds a constant 5/6 padding offset) ---
$ws.Range("G1").ColumnWidth = 17 - (5/6)
